$wb = $excel.ActiveWorkbook

# --- Table S2 - Plasticity AIC: remove AIC_wt (col E) and BIC_wt (col G) ---
$ws2 = $wb.Worksheets.Item("Table S2 - Plasticity AIC")
$ws2.Range("G:G").Delete()
$ws2.Range("E:E").Delete()
$ws2.Range("E:E").ColumnWidth = 4.83
$ws2.Range("F:F").ColumnWidth = 13.83
$ws2.Range("G:G").ColumnWidth = 10.83

# --- Table S7 - HostVsymb Plast AIC: remove AIC_wt (col D) and BIC_wt (col F) ---
$ws7 = $wb.Worksheets.Item("Table S7 - HostVsymb Plast AIC")
$ws7.Range("F:F").Delete()
$ws7.Range("D:D").Delete()
$ws7.Range("D:D").ColumnWidth = 5.83
$ws7.Range("E:E").ColumnWidth = 13.83
$ws7.Range("F:F").ColumnWidth = 10.83
$ws7.Range("G:G").ColumnWidth = 10.83
